$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as text in the source data (e.g. "26.381.59",
# "0.000007210"); Excel auto-converts numeric-looking strings and trims
# meaningful trailing/formatting digits, so force text format first, then
# clear the number-format override back to the default "Normal" style so no
# stray style index is left on the cell.
$dAddrs = @('D2', 'D3', 'D4', 'D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $dAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.381.59'
$ws.Range('D3').Value = '1.723.80'
$ws.Range('D4').Value = '0.9992'
$ws.Range('D5').Value = '242.40'
$ws.Range('D6').Value = '0.9996'
$ws.Range('D8').Value = '0.2590'
$ws.Range('D9').Value = '0.06190'
$ws.Range('D10').Value = '1.722.89'
$ws.Range('D11').Value = '0.06983'
$ws.Range('D12').Value = '15.52'
$ws.Range('D13').Value = '4.524'
$ws.Range('D14').Value = '0.5971'
$ws.Range('D15').Value = '77.18'
$ws.Range('D16').Value = '0.9993'
$ws.Range('D17').Value = '26.390.14'
$ws.Range('D18').Value = '0.9994'
$ws.Range('D19').Value = '0.000007210'
$ws.Range('D20').Value = '11.31'
$ws.Range('D21').Value = '1.948.93'
$ws.Range('D22').Value = '4.443'
$ws.Range('D24').Value = '5.091'
$ws.Range('D25').Value = '137.87'
$ws.Range('D26').Value = '15.23'
$ws.Range('D27').Value = '1.400'
$ws.Range('D29').Value = '1.725'
$ws.Range('D30').Value = '3.903'
$ws.Range('D31').Value = '0.08011'
$ws.Range('D35').Value = '0.9985'
$ws.Range('D36').Value = '0.6228'
$ws.Range('D37').Value = '0.9279'
$ws.Range('D38').Value = '1.959'
$ws.Range('D39').Value = '2.384'
$ws.Range('D40').Value = '0.9986'
$ws.Range('D41').Value = '100.84'
$ws.Range('D42').Value = '0.01475'
$ws.Range('D43').Value = '5.471'
$ws.Range('D44').Value = '0.3840'
$ws.Range('D45').Value = '6.907'
$ws.Range('D46').Value = '0.1162'
$ws.Range('D47').Value = '0.05371'
$ws.Range('D48').Value = '30.16'
$ws.Range('D49').Value = '7.667'
$ws.Range('D50').Value = '1.223'
$ws.Range('D51').Value = '51.03'

foreach ($addr in $dAddrs) {
    $ws.Range($addr).Style = "Normal"
}

# Columns B, C, E are plain text (names / URLs / padded percentages) and
# are not at risk of numeric auto-conversion.
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E5').Value = '  -1.91%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -2.95%  '
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('E13').Value = '  -2.91%  '
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('E24').Value = '  -3.15%  '
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  -1.58%  '
$ws.Range('E30').Value = '  -1.72%  '
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('E33').Value = '  -1.79%  '
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('E36').Value = '  -2.07%  '
$ws.Range('E37').Value = '  +3.22%  '
$ws.Range('E38').Value = '  -2.96%  '
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E42').Value = '  -2.23%  '
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('E44').Value = '  -1.15%  '
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('E46').Value = '  -1.78%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('E49').Value = '  -1.89%  '
$ws.Range('E50').Value = '  -2.04%  '
$ws.Range('E51').Value = '  -0.63%  '
